$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary rows: MEDIA / DESVIO PADRAO / TAXA DE SUCESSO ---------------

$ws.Range("A102").Value = "MEDIA"
$ws.Range("A103").Value = "DESVIO PADRAO"
$ws.Range("A104").Value = "TAXA DE SUCESSO"

$ws.Range("B102").Formula = "=MEDIAN(B2:B101)"
$ws.Range("C102").Formula = "=MEDIAN(C2:C101)"
$ws.Range("D102").Formula = "=MEDIAN(D2:D101)"

$ws.Range("B103").Formula = "=STDEV.P(B2:B101)"
$ws.Range("C103").Formula = "=STDEV.P(C2:C101)"
$ws.Range("D103").Formula = "=STDEV.P(D2:D101)"

$ws.Range("B104").Formula = "=COUNTIF(B2:B101,"">=4"")"
$ws.Range("C104").Formula = "=COUNTIF(C2:C101,"">=4"")"
$ws.Range("D104").Formula = "=COUNTIF(D2:D101,"">=4"")"

# --- Formatting ---------------------------------------------------------------

# Label cells (A102:A104): bold white text on a black fill, left/center aligned.
# Built up on an off-grid helper cell first and pasted across so the three
# label cells all land on a single shared style entry.
$tmpl = $ws.Range("Z1")
$tmpl.Interior.ThemeColor = 1
$tmpl.Font.Bold = $true
$tmpl.Font.ColorIndex = 2
$tmpl.HorizontalAlignment = -4131
$tmpl.VerticalAlignment = -4108
$tmpl.Copy()
$ws.Range("A102:A104").PasteSpecial(-4122)
$tmpl.Clear()

# Numeric cells (B102:D104): built-in "Comma" number style.
$ws.Range("B102:D104").Style = "Comma"

# --- Column widths -------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 16.16666666666667
$ws.Columns.Item(8).ColumnWidth = 5471.666666666667

# --- View state ------------------------------------------------------------

$null = $ws.Range("B104:D104").Select()
